$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31..103 down to 32..104
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with its data
$ws.Cells.Item(31, 1).Value = 5
$ws.Cells.Item(31, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(31, 3).Value = "Maule"
$ws.Cells.Item(31, 4).Value = 44935
$ws.Cells.Item(31, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31, 5).Value = 7
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100101
$ws.Cells.Item(31, 8).Value = "Berries"
$ws.Cells.Item(31, 9).Value = 100101001
$ws.Cells.Item(31, 10).Value = "Arándano (blue)"
$ws.Cells.Item(31, 11).Value = "Sin especificar"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 40
$ws.Cells.Item(31, 14).Value = 3000
$ws.Cells.Item(31, 15).Value = 3000
$ws.Cells.Item(31, 16).Value = 3000
$ws.Cells.Item(31, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(31, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(31, 19).Value = 1500
$ws.Cells.Item(31, 20).Value = 2
